$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 (shifts existing rows 2-4 down to rows 3-5)
$ws.Rows("2:2").Insert()

# The inserted row picks up the header row's formatting by default - reset it
# back to plain/default formatting like the other data rows.
$ws.Range("A2:R2").ClearFormats()

# Column D (Fecha) uses the date/time number-format style shared by the other
# data rows - match it.
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Populate the new weekly price record (same market/category, new date & prices)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 45092
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100114002
$ws.Range("G2").Value = "Camote"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 210
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 10714
$ws.Range("N2").Value = "$/malla 18 kilos"
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 595
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = "Hortaliza"
